$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get a plain numeric-looking string (e.g. "115.57",
# "0.0932") that Excel would otherwise auto-convert to a number (losing
# trailing zeros / switching to scientific notation). Force Text format
# on just those so the exact string round-trips.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.894.58'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '2.242.90'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '115.57'
$ws.Range("E5").Value = '  +2.60%  '
$ws.Range("D6").Value = '287.10'
$ws.Range("E6").Value = '  +8.19%  '
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  -3.62%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.615'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '46.80'
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").Value = '0.0932'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").Value = '9.18'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("D14").Value = '15.44'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").Value = '0.886'
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("D16").Value = '2.578.48'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").Value = '2.242.60'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").Value = '42.803.53'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("D20").Value = '6.99'
$ws.Range("E20").Value = '  +3.63%  '
$ws.Range("D21").Value = '73.64'
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").Value = '3.28'
$ws.Range("E22").Value = '  +13.41%  '
$ws.Range("D23").Value = '2.37'
$ws.Range("E23").Value = '  -2.58%  '
$ws.Range("D24").Value = '232.49'
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").Value = '9.22'
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("D26").Value = '12.15'
$ws.Range("E26").Value = '  +5.67%  '
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("D28").Value = '3.95'
$ws.Range("E28").Value = '  -0.68%  '
$ws.Range("D29").Value = '40.22'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.24'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("B31").Value = 'WEMIXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '3.29'
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").Value = '175.49'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").Value = '21.22'
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").Value = '0.0908'
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("D35").Value = '4.60'
$ws.Range("E35").Value = '  +18.71%  '
$ws.Range("D36").Value = '5.62'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  -2.94%  '
$ws.Range("D38").Value = '0.0374'
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").Value = '4.64'
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("D41").Value = '2.66'
$ws.Range("E41").Value = '  +2.17%  '
$ws.Range("D42").Value = '73.06'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").Value = '13.60'
$ws.Range("E43").Value = '  -4.93%  '
$ws.Range("D44").Value = '0.235'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = '1.34'
$ws.Range("E46").Value = '  -1.94%  '
$ws.Range("D47").Value = '5.62'
$ws.Range("E47").Value = '  -7.30%  '
$ws.Range("E48").Value = '  +2.79%  '
$ws.Range("D49").Value = '8.57'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").Value = '0.652'
$ws.Range("E50").Value = '  +5.60%  '
$ws.Range("D51").Value = '0.474'
$ws.Range("E51").Value = '  +8.41%  '
